$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.356401443481445
$ws.Range("B1").Value = 5.543663024902344
$ws.Range("C1").Value = 3.67113733291626
$ws.Range("D1").Value = 0.9658727645874023
$ws.Range("E1").Value = 0.6174454092979431
